$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Best-effort: nudge the saved window position to match the author's
# recorded xWindow (cosmetic, screen-position-only attribute).
try {
    $win = $excel.Workbooks.Item(1).Windows.Item(1)
    $win.Left = 15740
} catch {
    # Not fatal if the host doesn't track this.
}

# New "featureRequiringCoverage" column (F) added to the drug table,
# mapping each drug's HCV protein category to its short code (NS3,
# NS5A, NS5B) for subgenomic-sequence coverage checks.
# Data cells are written first, header last, so new shared strings are
# appended in the same order as the target workbook (NS3, NS5A, NS5B,
# then the header text).
$ws.Range("F2").Value = "NS3"
$ws.Range("F3").Value = "NS5A"
$ws.Range("F4").Value = "NS5A"
$ws.Range("F5").Value = "NS3"
$ws.Range("F6").Value = "NS5B"
$ws.Range("F1").Value = "featureRequiringCoverage"

# Match the author's final selection: a single active cell at F1
# (previously the whole sheet, A1:XFD1048576, was selected).
$ws.Range("F1").Select()
